# Actualización automática 2025-06-03 09:20:07
# Update the "PRESUPUESTO" (budget) column (G) values on the "VENTA MENSUAL"
# sheet for a set of clients, and refresh the total in the last row (G55).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$updates = @{
    2  = 1500
    5  = 4000
    6  = 2000
    11 = 2000
    12 = 3000
    19 = 3000
    22 = 2000
    24 = 2000
    26 = 15000
    27 = 3000
    29 = 5000
    31 = 3000
    36 = 10000
    37 = 2000
    39 = 2000
    41 = 2000
    42 = 3000
    43 = 4000
    44 = 3000
    45 = 1500
    47 = 3000
    48 = 3000
    49 = 2000
    50 = 5000
    52 = 2000
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}

# Update the total row with the new sum of the PRESUPUESTO column.
$ws.Range("G55").Value = 88000
